$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps on the existing "data" sheet ---
$ws.Range("F2").Value = "2021-10-05 14:20:31.497042"
$ws.Range("F3").Value = "2021-10-05 14:20:31.497050"
$ws.Range("F4").Value = "2021-10-05 14:20:31.497053"
$ws.Range("F5").Value = "2021-10-05 14:20:31.497056"
$ws.Range("F6").Value = "2021-10-05 14:20:31.497059"
$ws.Range("F7").Value = "2021-10-05 14:20:31.497062"
$ws.Range("F8").Value = "2021-10-05 14:20:31.497064"
$ws.Range("F9").Value = "2021-10-05 14:20:31.497067"
$ws.Range("F10").Value = "2021-10-05 14:20:31.497070"
$ws.Range("F11").Value = "2021-10-05 14:20:31.497072"
$ws.Range("F12").Value = "2021-10-05 14:20:31.497075"
$ws.Range("F13").Value = "2021-10-05 14:20:31.497078"
$ws.Range("F14").Value = "2021-10-05 14:20:31.497080"

# --- Add a new "metadata" sheet right after the "data" sheet ---
$newSheet = $wb.Worksheets.Add($null, $ws)
$newSheet.Name = "metadata"

# Header row (B1:G1)
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row (A2:G2)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Genodermatoses with malignancies"
$newSheet.Range("C2").Value = 201

# "data_version" must stay textual ("1.6"), not become a float -> force text
# format first, then fix the resulting cell style back to the default below.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.6"

$newSheet.Range("E2").Value = "2019-10-09T08:30:18.879174Z"
$newSheet.Range("F2").Value = "2021-10-05 14:20:31.493605"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/201/?format=json"

# --- Match formatting of the "data" sheet's header/index-column styles ---

# Header row B1:G1 -> same bold/centered/bordered style as "data"!B1
$ws.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)

# Index cell A2 -> same style as "data"!A2 (the index-column style)
$ws.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# D2 picked up a stray text-number-format style above; reset it back to the
# default (unstyled) look, matching every other plain data cell on the row.
$ws.Range("B2").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)

$newSheet.Range("A1").Select() | Out-Null
